$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Every Price/Volume cell on this sheet is stored as plain text (inline
# strings in the source file), never as a Number. Values that look numeric are
# written with a leading apostrophe (quote-prefix) so Excel keeps them as text
# exactly as given -- otherwise Excel would auto-convert them to a Number and
# silently drop meaningful trailing zeros (e.g. "0.890" -> 0.89).

# --- Rows 44 and 45: ranking swap (Aave moves up to #44, Cronos drops to #45) ---
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'90.03"
$ws.Range("E44").Value = "  -5.24%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0877"
$ws.Range("E45").Value = "  -7.23%  "

# --- Remaining per-cell price / volume updates ---
$ws.Range("D2").Value = "35.329.43"
$ws.Range("E2").Value = "  -3.75%  "
$ws.Range("D3").Value = "1.977.24"
$ws.Range("E3").Value = "  -4.86%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'239.63"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("E6").Value = "  -11.80%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'55.93"
$ws.Range("E8").Value = "  +6.68%  "
$ws.Range("D9").Value = "'59.18"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("E11").Value = "  -4.42%  "
$ws.Range("D12").Value = "'0.103"
$ws.Range("E12").Value = "  -6.32%  "
$ws.Range("D13").Value = "'0.890"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "'14.16"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("D15").Value = "2.268.77"
$ws.Range("E15").Value = "  -4.88%  "
$ws.Range("E16").Value = "  -3.96%  "
$ws.Range("D17").Value = "1.976.74"
$ws.Range("E17").Value = "  -5.27%  "
$ws.Range("D18").Value = "'16.95"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").Value = "35.237.79"
$ws.Range("E19").Value = "  -3.89%  "
$ws.Range("D20").Value = "'69.65"
$ws.Range("E20").Value = "  -3.85%  "
$ws.Range("E21").Value = "  -4.41%  "
$ws.Range("D22").Value = "'230.95"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("E23").Value = "  -7.14%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  -5.54%  "
$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "  +4.32%  "
$ws.Range("D27").Value = "'162.60"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  -5.41%  "
$ws.Range("D29").Value = "'19.34"
$ws.Range("E29").Value = "  -5.62%  "
$ws.Range("D30").Value = "'0.117"
$ws.Range("E30").Value = "  -10.01%  "
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").Value = "'4.73"
$ws.Range("E32").Value = "  -8.00%  "
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("D34").Value = "'0.0893"
$ws.Range("E34").Value = "  +9.38%  "
$ws.Range("E35").Value = "  -9.10%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -5.59%  "
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").Value = "'4.81"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("E40").Value = "  -5.91%  "
$ws.Range("D41").Value = "'2.80"
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("E42").Value = "  -5.40%  "
$ws.Range("E43").Value = "  -5.97%  "
$ws.Range("D46").Value = "'7.41"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").Value = "1.341.82"
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("D48").Value = "'15.30"
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").Value = "'2.85"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "'2.23"
$ws.Range("E50").Value = "  -6.34%  "
$ws.Range("D51").Value = "'45.13"
$ws.Range("E51").Value = "  -0.39%  "
